# Simplified logic for top-node processing
# The content placeholder's bottom paragraphs ("Bold", "Item one", "Item
# two", "Item three") were being rendered as a bulleted list even though
# the first two paragraphs ("Both"/"Hello", "Italic"/"World") already
# suppress their bullet. Bring the remaining paragraphs in line (no
# bullet) and restore the intended emphasis (underline on "Item one" and
# "Item two", plus bold on "Item two").

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item("Content Placeholder 2")
$tr = $shape.TextFrame.TextRange

$paraBold = $tr.Paragraphs(3)
$paraBold.ParagraphFormat.Bullet.Visible = $false

$paraItemOne = $tr.Paragraphs(4)
$paraItemOne.ParagraphFormat.Bullet.Visible = $false
$paraItemOne.Font.Underline = $true

$paraItemTwo = $tr.Paragraphs(5)
$paraItemTwo.ParagraphFormat.Bullet.Visible = $false
$paraItemTwo.Font.Bold = $true
$paraItemTwo.Font.Underline = $true

$paraItemThree = $tr.Paragraphs(6)
$paraItemThree.ParagraphFormat.Bullet.Visible = $false
